# Auto-generated: updates numeric market-data cells (columns H-N) across all
# class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the latest
# scheduled-runner price snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 782.8333
$ws.Range("I32").Value = 583.8889
$ws.Range("K32").Value = 583.8889
$ws.Range("M32").Value = -257.8889
$ws.Range("H98").Value = 6144.4814
$ws.Range("I98").Value = 4018.6667
$ws.Range("J98").Value = 8801.75
$ws.Range("K98").Value = 4018.6667
$ws.Range("L98").Value = 8801.75
$ws.Range("M98").Value = -2520.6667
$ws.Range("N98").Value = -11797.75
$ws.Range("H122").Value = 6144.4814
$ws.Range("I122").Value = 4018.6667
$ws.Range("J122").Value = 8801.75
$ws.Range("K122").Value = 12056.0001
$ws.Range("L122").Value = 26405.25
$ws.Range("M122").Value = -9606.000100000001
$ws.Range("N122").Value = -31305.25
$ws.Range("H129").Value = 2611.5908
$ws.Range("I129").Value = 478.5
$ws.Range("J129").Value = 2824.9
$ws.Range("K129").Value = 1435.5
$ws.Range("L129").Value = 8474.700000000001
$ws.Range("M129").Value = 3564.5
$ws.Range("N129").Value = -18474.7
$ws.Range("H137").Value = 1137682.4
$ws.Range("I137").Value = 4768902
$ws.Range("J137").Value = 2926.1875
$ws.Range("K137").Value = 14306706
$ws.Range("L137").Value = 8778.5625
$ws.Range("M137").Value = -14304156
$ws.Range("N137").Value = -13878.5625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3581.5
$ws.Range("I45").Value = 2747.375
$ws.Range("J45").Value = 5249.75
$ws.Range("K45").Value = 2747.375
$ws.Range("L45").Value = 5249.75
$ws.Range("M45").Value = -2370.375
$ws.Range("N45").Value = -6003.75
$ws.Range("H74").Value = 274517.5
$ws.Range("I74").Value = 627727.4
$ws.Range("J74").Value = 1582.6364
$ws.Range("K74").Value = 627727.4
$ws.Range("L74").Value = 1582.6364
$ws.Range("M74").Value = -626853.4
$ws.Range("N74").Value = -3330.6364
$ws.Range("H77").Value = 274517.5
$ws.Range("I77").Value = 627727.4
$ws.Range("J77").Value = 1582.6364
$ws.Range("K77").Value = 3138637
$ws.Range("L77").Value = 7913.182000000001
$ws.Range("M77").Value = -3134269
$ws.Range("N77").Value = -16649.182

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2553.2856
$ws.Range("I86").Value = 1943
$ws.Range("J86").Value = 2797.4
$ws.Range("K86").Value = 1943
$ws.Range("L86").Value = 2797.4
$ws.Range("M86").Value = -820
$ws.Range("N86").Value = -5043.4
$ws.Range("H89").Value = 2553.2856
$ws.Range("I89").Value = 1943
$ws.Range("J89").Value = 2797.4
$ws.Range("K89").Value = 9715
$ws.Range("L89").Value = 13987
$ws.Range("M89").Value = -4099
$ws.Range("N89").Value = -25219
$ws.Range("H94").Value = 654.2857
$ws.Range("I94").Value = 555
$ws.Range("K94").Value = 555
$ws.Range("M94").Value = -104
$ws.Range("H99").Value = 1238.3077
$ws.Range("I99").Value = 1116.5
$ws.Range("J99").Value = 2700
$ws.Range("K99").Value = 1116.5
$ws.Range("L99").Value = 2700
$ws.Range("M99").Value = 381.5
$ws.Range("N99").Value = -5696

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2027.5209
$ws.Range("I58").Value = 1167.875
$ws.Range("J58").Value = 3746.8125
$ws.Range("K58").Value = 1167.875
$ws.Range("L58").Value = 3746.8125
$ws.Range("M58").Value = -964.875
$ws.Range("N58").Value = -4152.8125
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882
$ws.Range("H132").Value = 3550.138
$ws.Range("I132").Value = 2984.9473
$ws.Range("J132").Value = 4624
$ws.Range("K132").Value = 8954.841899999999
$ws.Range("L132").Value = 13872
$ws.Range("M132").Value = -6424.841899999999
$ws.Range("N132").Value = -18932
$ws.Range("H136").Value = 2027.5209
$ws.Range("I136").Value = 1167.875
$ws.Range("J136").Value = 3746.8125
$ws.Range("K136").Value = 3503.625
$ws.Range("L136").Value = 11240.4375
$ws.Range("M136").Value = -953.625
$ws.Range("N136").Value = -16340.4375
$ws.Range("H137").Value = 46820
$ws.Range("J137").Value = 46820
$ws.Range("L137").Value = 46820
$ws.Range("N137").Value = -57020

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2747.1177
$ws.Range("I5").Value = 1536.5454
$ws.Range("J5").Value = 4966.5
$ws.Range("K5").Value = 4609.6362
$ws.Range("L5").Value = 14899.5
$ws.Range("M5").Value = -4497.6362
$ws.Range("N5").Value = -15123.5
$ws.Range("H68").Value = 2943.2952
$ws.Range("I68").Value = 1409
$ws.Range("J68").Value = 3208.8462
$ws.Range("K68").Value = 4227
$ws.Range("L68").Value = 9626.5386
$ws.Range("M68").Value = -3416
$ws.Range("N68").Value = -11248.5386
$ws.Range("H71").Value = 2943.2952
$ws.Range("I71").Value = 1409
$ws.Range("J71").Value = 3208.8462
$ws.Range("K71").Value = 12681
$ws.Range("L71").Value = 28879.6158
$ws.Range("M71").Value = -8625
$ws.Range("N71").Value = -36991.6158
$ws.Range("H92").Value = 1532.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1532.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 4597.5
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -7093.5
$ws.Range("H113").Value = 455.67307
$ws.Range("I113").Value = 451.79413
$ws.Range("J113").Value = 463
$ws.Range("K113").Value = 1355.38239
$ws.Range("L113").Value = 1389
$ws.Range("M113").Value = 814.61761
$ws.Range("N113").Value = -5729
$ws.Range("H131").Value = 789.87
$ws.Range("I131").Value = 474
$ws.Range("J131").Value = 817.337
$ws.Range("K131").Value = 1422
$ws.Range("L131").Value = 2452.011
$ws.Range("M131").Value = 3618
$ws.Range("N131").Value = -12532.011
$ws.Range("H132").Value = 2352.5789
$ws.Range("I132").Value = 990
$ws.Range("J132").Value = 2715.9333
$ws.Range("K132").Value = 8910
$ws.Range("L132").Value = 24443.3997
$ws.Range("M132").Value = -6380
$ws.Range("N132").Value = -29503.3997
$ws.Range("H135").Value = 2747.1177
$ws.Range("I135").Value = 1536.5454
$ws.Range("J135").Value = 4966.5
$ws.Range("K135").Value = 13828.9086
$ws.Range("L135").Value = 44698.5
$ws.Range("M135").Value = -11293.9086
$ws.Range("N135").Value = -49768.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 32443
$ws.Range("J111").Value = 32443
$ws.Range("L111").Value = 32443
$ws.Range("N111").Value = -38577
$ws.Range("H123").Value = 10497.6
$ws.Range("J123").Value = 10497.6
$ws.Range("L123").Value = 10497.6
$ws.Range("N123").Value = -15397.6
$ws.Range("H132").Value = 4535.6484
$ws.Range("I132").Value = 3895.1667
$ws.Range("J132").Value = 5142.421
$ws.Range("K132").Value = 11685.5001
$ws.Range("L132").Value = 15427.263
$ws.Range("M132").Value = -9155.500100000001
$ws.Range("N132").Value = -20487.263

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7409676.5
$ws.Range("I93").Value = 8548857
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 8548857
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -8547609
$ws.Range("N93").Value = -7496
$ws.Range("H100").Value = 5033.6665
$ws.Range("I100").Value = 2075.75
$ws.Range("J100").Value = 7400
$ws.Range("K100").Value = 2075.75
$ws.Range("L100").Value = 7400
$ws.Range("M100").Value = -1534.75
$ws.Range("N100").Value = -8482
$ws.Range("H122").Value = 4064.6667
$ws.Range("I122").Value = 3208
$ws.Range("J122").Value = 5521
$ws.Range("K122").Value = 9624
$ws.Range("L122").Value = 16563
$ws.Range("M122").Value = -7174
$ws.Range("N122").Value = -21463

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 78755500
$ws.Range("I96").Value = 90909910
$ws.Range("J96").Value = 11906262
$ws.Range("K96").Value = 90909910
$ws.Range("L96").Value = 11906262
$ws.Range("M96").Value = -90908537
$ws.Range("N96").Value = -11909008
$ws.Range("H122").Value = 3925.5789
$ws.Range("I122").Value = 1022.3333
$ws.Range("J122").Value = 6538.5
$ws.Range("K122").Value = 3066.9999
$ws.Range("L122").Value = 19615.5
$ws.Range("M122").Value = -616.9998999999998
$ws.Range("N122").Value = -24515.5
$ws.Range("H126").Value = 712159.75
$ws.Range("I126").Value = 2533.3333
$ws.Range("J126").Value = 1185244
$ws.Range("K126").Value = 7599.999899999999
$ws.Range("L126").Value = 3555732
$ws.Range("M126").Value = -5129.999899999999
$ws.Range("N126").Value = -3560672

